$wb = $excel.ActiveWorkbook

# --- Worksheet 1: "RUNMANAGER" (testcasename/execute/priority/count table) ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new "description" column after column A (shifts execute/priority/count right)
$ws1.Columns.Item(2).Insert() | Out-Null

# Set width for the new description column (best effort; engine rounds to pixel grid)
$ws1.Columns.Item(2).ColumnWidth = 51.42578125

# Header row
$ws1.Range("B1").Value = "description"

# Data rows - new description text in column B
$ws1.Range("B2").Value = "To verify if user can successfully login and logout"
$ws1.Range("B3").Value = "To verify if error is prompted on invalid login"
$ws1.Range("B4").Value = "To verify user is able to access Admin page"
$ws1.Range("B5").Value = "To verify footer is present"
$ws1.Range("B6").Value = "'"

# Priority/count columns (now D & E after insert) - stored as text (quote-prefixed numbers)
$ws1.Range("D2").Value = "'1"
$ws1.Range("E2").Value = "'2"
$ws1.Range("D3").Value = "'2"
$ws1.Range("E3").Value = "'1"
$ws1.Range("D4").Value = "'3"
$ws1.Range("E4").Value = "'1"
$ws1.Range("D5").Value = "'4"
$ws1.Range("D6").Value = "'"

# --- Worksheet 2: "testData" (testcasename/username/password/assertion value/execute table) ---
$ws2 = $wb.Worksheets.Item(2)

# Rename sheet testData -> TESTDATA
$ws2.Name = "TESTDATA"

# Append new row 8, copying the formatting (borders) of the row above first
$ws2.Range("A7:E7").Copy() | Out-Null
$ws2.Range("A8:E8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws2.Range("A8").Value = "loginLogoutTest"
$ws2.Range("B8").Value = "Admin123"
$ws2.Range("C8").Value = "admin123"
$ws2.Range("D8").Value = "OrangeHRM"
$ws2.Range("E8").Value = "Yes"

# --- Selections: set sheet1 selection first (so it is not left as the active tab), then sheet2 last (active tab) ---
$ws1.Range("E2").Select() | Out-Null
$ws2.Range("A8").Select() | Out-Null
